# "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The "Periodo Mora" (E16:E25) and "Valor Mora" (F16:F25) columns are
# reshuffled: the old period list (1712, 1801..1809, oldest first) is
# replaced by the new list in the opposite order (1809, 1808 .. 1801,
# 1712, newest first), carrying its matching arrears value along with
# it. "Salario Basico" (column G) is identical for every period so it
# is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 16
$lastRow = 25

# Snapshot the current Periodo Mora / Valor Mora values before
# overwriting anything. (Value2 is used for the read-back since it
# reflects literals cleanly; Value is still used for the write below.)
$periodos = @()
$valores = @()
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $periodos += $ws.Cells.Item($r, 5).Value2
    $valores += $ws.Cells.Item($r, 6).Value2
}

# Write them back in reverse row order.
$rowCount = $lastRow - $firstRow + 1
for ($i = 0; $i -lt $rowCount; $i++) {
    $targetRow = $firstRow + $i
    $sourceIndex = $rowCount - 1 - $i
    $ws.Cells.Item($targetRow, 5).Value = $periodos[$sourceIndex]
    $ws.Cells.Item($targetRow, 6).Value = $valores[$sourceIndex]
}
